$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the pre-existing "_GoBack" bookmark that split the cached table
#    number field result "17" into "1" + "7" (Table 17.--  Mean Number CPUE
#    ...), and re-merge the digits back into a single run. We locate the
#    digits by scanning backward from the unique "Mean Number CPUE" anchor
#    text, skipping the field/". -- " boilerplate in between.
# ---------------------------------------------------------------------------
$anchor = $d.Content.Duplicate
$gotAnchor = $anchor.Find.Execute("Mean Number CPUE")
if ($gotAnchor) {
    $mstart = $anchor.Start
    $scanBack = $d.Range($mstart - 25, $mstart)
    $n = $scanBack.Characters.Count
    $digitEndOffset = -1
    $digitStartOffset = -1
    for ($i = $n - 1; $i -ge 0; $i--) {
        $c = $scanBack.Characters.Item($i + 1)
        $t = $c.Text
        if ($t -match '^[0-9]$') {
            if ($digitEndOffset -eq -1) {
                $digitEndOffset = $c.Start + 1
            }
            $digitStartOffset = $c.Start
        } elseif ($digitEndOffset -ne -1) {
            break
        }
    }

    if ($digitStartOffset -ge 0) {
        $digits = $d.Range($digitStartOffset, $digitEndOffset)
        $digitsText = $digits.Text
        $d.Bookmarks("_GoBack").Delete()
        $digits.Find.Execute($digitsText, $false, $false, $false, $false, `
            $false, $true, 1, $false, $digitsText, 2) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# 2) Heading3 "H3 Selected Invertebrates Estimates" -> split into runs with a
#    new "_GoBack" bookmark wrapped around the "I" of "Invertebrates" (mirrors
#    the editing artifact Word leaves behind after the cursor was last
#    positioned there).
# ---------------------------------------------------------------------------
$findRng = $d.Content.Duplicate
$found = $findRng.Find.Execute("H3 Selected Invertebrates Estimates")
if ($found) {
    $base = $findRng.Start

    # Split boundary right before "I" (after "H3 Selected ") - force a run
    # break with a throw-away bookmark that we immediately delete.
    $p1 = $d.Range($base + 12, $base + 12)
    $d.Bookmarks.Add("ZZZTempSplit1", $p1) | Out-Null
    $d.Bookmarks("ZZZTempSplit1").Delete()

    # Wrap the single "I" character in the real "_GoBack" bookmark.
    $p2 = $d.Range($base + 13, $base + 13)
    $d.Bookmarks.Add("_GoBack", $p2) | Out-Null

    # Split boundary right after "nvertebrates" (before " Estimates").
    $p3 = $d.Range($base + 25, $base + 25)
    $d.Bookmarks.Add("ZZZTempSplit2", $p3) | Out-Null
    $d.Bookmarks("ZZZTempSplit2").Delete()
}

# ---------------------------------------------------------------------------
# 3) Heading styles: bump "space after" from 0 to 12pt (240 twentieths) for
#    Heading 1/2/3.
# ---------------------------------------------------------------------------
$d.Styles("Heading 1").ParagraphFormat.SpaceAfter = 12
$d.Styles("Heading 2").ParagraphFormat.SpaceAfter = 12
$d.Styles("Heading 3").ParagraphFormat.SpaceAfter = 12

# ---------------------------------------------------------------------------
# 4) Materialize word/footnotes.xml + word/endnotes.xml (the document
#    currently has neither part). Adding then immediately deleting a
#    footnote mints both parts with Word's standard separator/
#    continuation-separator boilerplate, matching a google-docs/drive
#    round-trip import.
# ---------------------------------------------------------------------------
$fnRange = $d.Range(0, 0)
$fn = $d.Footnotes.Add($fnRange)
$fn.Delete()

Write-Output "done"
